# Updated symbol list with refreshed cryptocurrency price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        $Sheet,
        [string]$CellRef,
        [string]$Text
    )
    $range = $Sheet.Range($CellRef)
    # Force text storage so numeric-looking strings (prices, percentages)
    # are NOT reinterpreted as numbers/percentages by Excel's input parser.
    $range.NumberFormat = "@"
    $range.Value = $Text
    # Drop the temporary "Text" number format again so the cell's
    # formatting matches the original (unformatted) cell.
    $range.ClearFormats()
}

Set-CellText $ws "D2" "244.97"
Set-CellText $ws "E2" "-0.67%"
Set-CellText $ws "D3" "28.63"
Set-CellText $ws "E3" "-3.99%"
Set-CellText $ws "D4" "5.236"
Set-CellText $ws "E4" "0.96%"
Set-CellText $ws "D5" "0.05699"
Set-CellText $ws "D6" "6.616"
Set-CellText $ws "E6" "0.27%"
Set-CellText $ws "D7" "3.188"
Set-CellText $ws "E7" "3.25%"
Set-CellText $ws "D8" "0.8508"
Set-CellText $ws "E8" "-0.70%"
Set-CellText $ws "D9" "0.8540"
Set-CellText $ws "E9" "-1.89%"
Set-CellText $ws "E10" "0.59%"
Set-CellText $ws "D11" "0.07091"
Set-CellText $ws "E11" "0.73%"
Set-CellText $ws "D12" "0.03155"
Set-CellText $ws "E12" "7.90%"
Set-CellText $ws "D13" "0.09202"
Set-CellText $ws "D14" "0.001541"
Set-CellText $ws "E14" "1.76%"
Set-CellText $ws "D15" "0.0005969"
Set-CellText $ws "E15" "-0.85%"
Set-CellText $ws "D16" "0.005994"
Set-CellText $ws "E16" "-1.63%"
Set-CellText $ws "D17" "3.493"
Set-CellText $ws "E17" "0.20%"
Set-CellText $ws "E18" "-4.45%"
Set-CellText $ws "E19" "0.36%"
Set-CellText $ws "D20" "0.03264"
Set-CellText $ws "E20" "-3.15%"
Set-CellText $ws "D21" "0.1298"
Set-CellText $ws "E21" "-1.34%"
Set-CellText $ws "D22" "3.485"
Set-CellText $ws "E22" "0.61%"
Set-CellText $ws "E23" "-2.25%"
Set-CellText $ws "D25" "0.001222"
Set-CellText $ws "E25" "-0.02%"
Set-CellText $ws "E26" "-17.50%"
Set-CellText $ws "D27" "0.0001200"
Set-CellText $ws "E27" "-0.78%"
Set-CellText $ws "D28" "0.0001449"
Set-CellText $ws "D40" "0.03752"
Set-CellText $ws "E40" "0.10%"
Set-CellText $ws "D41" "0.1064"
Set-CellText $ws "E41" "-0.76%"
Set-CellText $ws "E42" "-35.31%"
Set-CellText $ws "D43" "0.002490"
Set-CellText $ws "E43" "24.55%"
Set-CellText $ws "D44" "0.009339"
Set-CellText $ws "E44" "-2.41%"
Set-CellText $ws "D45" "0.00005281"
Set-CellText $ws "E45" "1.06%"
Set-CellText $ws "E46" "0.00%"
Set-CellText $ws "D47" "0.07498"
Set-CellText $ws "E47" "15.90%"
Set-CellText $ws "E48" "-3.14%"
Set-CellText $ws "E49" "0.00%"
Set-CellText $ws "E50" "0.00%"
